$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Inscritos (E8) 38 -> 39
$ws.Range("E8").Value = 39

# Row 10: Inscritos (E10) 25 -> 26, Pagos (F10) 11 -> 12, Inscricoes homologadas (H10) 11 -> 12
$ws.Range("E10").Value = 26
$ws.Range("F10").Value = 12
$ws.Range("H10").Value = 12
